$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of H1 onto the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Add new column values (I and J) for each data row
$values = @{
    2  = @(9, 9)
    3  = @(9, 9)
    4  = @(8, 8)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(8, 8)
    8  = @(9, 9)
    9  = @(8, 8)
    10 = @(9, 9)
    11 = @(8, 9)
    12 = @(7, 8)
    13 = @(8, 9)
    14 = @(9, 9)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(7, 7)
    19 = @(6, 6)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
